$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planning effectif")

# "avancement du site" - log hours worked on various tasks (time-of-day
# formatted durations, stored as fraction-of-a-day numbers, e.g. 0.5h = 0.5/24)
$ws.Range("F8").Value  = 0.5 / 24
$ws.Range("G12").Value = 0.5 / 24
$ws.Range("G13").Value = 0.5 / 24
$ws.Range("G15").Value = 0.5 / 24
$ws.Range("F19").Value = 2 / 24
$ws.Range("F21").Value = 0.5 / 24
$ws.Range("G25").Value = 1 / 24
$ws.Range("G26").Value = 1.5 / 24
$ws.Range("G29").Value = 0.5 / 24
$ws.Range("G31").Value = 0.5 / 24
$ws.Range("G32").Value = 1 / 24

# "documentation" task row
$ws.Range("F38").Value = 2 / 24
$ws.Range("G38").Value = 2 / 24
